$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 0: remember & remove the existing last comment (currently id=13,
# "Se cambio a electas en 2019 porque no enviaron datos.") so that the
# numeric id it occupies is freed up for the new comment we are about
# to insert earlier in the document. Word will hand that id back out
# to the next Comments.Add call, and this comment will be re-created
# afterwards (picking up the next free id, i.e. 14), matching the
# renumbering 13->14 seen in the target diff.
# ------------------------------------------------------------------
$lastComment = $d.Comments.Item($d.Comments.Count)
$lastCommentText = $lastComment.Range.Text
$lastCommentAuthor = $lastComment.Author
$lastCommentInitial = $lastComment.Initial
$lastComment.Delete()

# ------------------------------------------------------------------
# Step 1: split the "Participación en los Consejos de Desarrollo por
# sexo, según cargo " run so that ", según cargo" gets struck through
# and receives a new comment explaining that it now lives in the
# expanded table.
# ------------------------------------------------------------------
$rStrike = $d.Content
$rStrike.Find.Execute(", según cargo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rStrike.Font.StrikeThrough = $true

$rComment = $d.Range($rStrike.Start, $rStrike.End + 1)
$newComment = $d.Comments.Add($rComment, "Según cargo está en información ampliada, ya que por espacio no cabía en una sola tabla.")
$newComment.Author = "Paula Natalia Galvez Molina"
$newComment.Initial = "PG"

# ------------------------------------------------------------------
# Step 2: re-create the comment that used to be id=13 on the same
# anchor ("en" inside "Mujeres electas para en alcaldías"). It will
# now be appended last, receiving the next free id (14).
# ------------------------------------------------------------------
$rAnchor = $d.Content
$rAnchor.Find.Execute("Mujeres electas para en alcaldías", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rEn = $rAnchor.Words.Item(4)
$rEn.MoveEnd(1, -1)

$restoredComment = $d.Comments.Add($rEn, $lastCommentText)
$restoredComment.Author = $lastCommentAuthor
$restoredComment.Initial = $lastCommentInitial
